$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column L ("Número de Relatório" report-number column)
$ws.Range("L1").Value = "Número de Relatório"

# Columns C (Part Number) and E (Semana) hold digit-only values that must be
# stored as TEXT (matching the rest of the sheet, which uses
# numberStoredAsText). Pre-format those two ranges as text before writing so
# Excel doesn't coerce them into numeric cells, then drop the temporary
# number-format mark so the cells keep their plain (unstyled) look.
$ws.Range("C14:C24").NumberFormat = "@"
$ws.Range("E14:E24").NumberFormat = "@"

# New report rows (14-24), columns C..L
$data = @(
    @("12345678", "KUKÃO LD", "20", "ENGENHARIA",    "matheus", "2º TURNO", "CMM GLOBAL", "ACOMPANHAMENTO", "PEÇA MUITO BOA.",       ""),
    @("12345678", "KUKÃO LD", "5",  "ENGENHARIA",    "luis",    "3º TURNO", "PAQUÍMETRO", "ACOMPANHAMENTO", "fghe 5t3y6 e5 y45y3e",  ""),
    @("12345678", "KUKÃO LD", "2",  "QUALIDADE",     "izaac",   "1º TURNO", "METRASCAN",  "ANÁLISE",        "r5ty6 y43 56y35",       ""),
    @("53490058", "KUKÃO LE", "20", "ENGENHARIA",    "matheus", "2º TURNO", "CMM GLOBAL", "ACOMPANHAMENTO", "PEÇA OK!",              ""),
    @("53490059", "KUKÃO LD", "20", "ENGENHARIA",    "matheus", "2º TURNO", "CMM GLOBAL", "ACOMPANHAMENTO", "peça ok!",              "C2025.0016"),
    @("53490059", "KUKÃO LD", "20", "ENGENHARIA",    "matheus", "2º TURNO", "CMM GLOBAL", "ACOMPANHAMENTO", "peça ok!",              "C2025.0017"),
    @("53490059", "KUKÃO LD", "20", "FERRAMENTARIA", "luis",    "3º TURNO", "CMM GLOBAL", "ANÁLISE",        "peça ok!",              "C2025.0018"),
    @("53490059", "KUKÃO LE", "20", "ENGENHARIA",    "izaac",   "ADM",      "CMM GLOBAL", "ANÁLISE",        "peça ok!",              "C2025.0019"),
    @("53490059", "KUKÃO LE", "20", "ENGENHARIA",    "luis",    "1º TURNO", "CMM GLOBAL", "TRYOUT",         "peça ok!",              "C2025.0020"),
    @("53490059", "KUKÃO LE", "20", "FERRAMENTARIA", "luis",    "1º TURNO", "CMM GLOBAL", "ACOMPANHAMENTO", "peça ok!",              "C2025.0021"),
    @("53490059", "KUKÃO LE", "20", "QUALIDADE",     "izaac",   "1º TURNO", "METRASCAN",  "GEOMETRIA",      "peça ok!",              "C2025.0022")
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 3; $col -le 12; $col++) {
        $val = $rowData[$col - 3]
        if ($val -ne "") {
            $ws.Cells.Item($row, $col).Value = $val
        }
    }
}
